$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename "Attribut..." to "Attribute..." for consistency ---
$ws.Range("D1").Value = "AttributeName"
$ws.Range("E1").Value = "AttributeName"
$ws.Range("F1").Value = "AttributeDescriptionIT"

# --- IfcSpace attribute rows (6-11) reshuffled ---
# Row 8: IsExternal -> PredefinedType
$ws.Range("D8").Value = "PredefinedType"
$ws.Range("E8").Value = "PredefinedType"

# Row 9: PredefinedType -> IsInteriorOrExteriorSpace, and it now carries the Pset_SpaceCommon marker
$ws.Range("D9").Value = "IsInteriorOrExteriorSpace"
$ws.Range("E9").Value = "IsInteriorOrExteriorSpace"
$ws.Range("C9").Value = "Pset_SpaceCommon"

# Row 10: IsInteriorOrExteriorSpace -> IsExternal, Pset_SpaceCommon marker removed
$ws.Range("D10").Value = "IsExternal"
$ws.Range("E10").Value = "IsExternal"
$ws.Range("C10").Clear()

# --- Column width adjustments ---
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(1).ColumnWidth()
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(7).ColumnWidth()
$ws.Range("J1:K1").EntireColumn.ColumnWidth = $ws.Columns.Item(12).ColumnWidth()
